# Update Name of Algo
# Corrects numeric values in the KNN imputation result data (Sheet1, columns A-E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -13.164
$ws.Range("B7").Value = 5.910000000000001
$ws.Range("A10").Value = -21.663
$ws.Range("E10").Value = 16.364
$ws.Range("A12").Value = -21.515
$ws.Range("E14").Value = 16.95
$ws.Range("B15").Value = 5.072000000000001
$ws.Range("A18").Value = -21.694
$ws.Range("C18").Value = -10.704
$ws.Range("C19").Value = -11.706
$ws.Range("B20").Value = 7.217000000000001
$ws.Range("C27").Value = -13.359
$ws.Range("B29").Value = 5.645
$ws.Range("B30").Value = 6.02
$ws.Range("B31").Value = 5.654000000000001
$ws.Range("E32").Value = 16.847
$ws.Range("E35").Value = 16.384
$ws.Range("A37").Value = -20.287
$ws.Range("B40").Value = 8.870999999999999
$ws.Range("C42").Value = -12.533
$ws.Range("E43").Value = 17.071
$ws.Range("C44").Value = -12.324
$ws.Range("C47").Value = -12.083
$ws.Range("E49").Value = 16.359
$ws.Range("A55").Value = -21.875
$ws.Range("E56").Value = 16.374
$ws.Range("C58").Value = -12.621
$ws.Range("A68").Value = -21.581
$ws.Range("B68").Value = 5.242
$ws.Range("E69").Value = 17.448
$ws.Range("C73").Value = -12.57
$ws.Range("B76").Value = 6.654000000000001
$ws.Range("A77").Value = -20.252
$ws.Range("A78").Value = -19.842
$ws.Range("E81").Value = 16.872
$ws.Range("B87").Value = 5.429
$ws.Range("B88").Value = 5.75
$ws.Range("E92").Value = 17.834
$ws.Range("C95").Value = -11.367
$ws.Range("B96").Value = 6.484
$ws.Range("B98").Value = 5.031999999999999
$ws.Range("B101").Value = 9.096
$ws.Range("C101").Value = -12.234
$ws.Range("B102").Value = 7.74
